$d = $word.ActiveDocument

# --- Part 1: Collapse the three CORE COMPETENCIES bullet paragraphs into one ---
$bullet = [char]0x2022

$coreFirst = $d.Paragraphs.Item(6)
$r = $coreFirst.Range
$r.MoveEnd(1, -1)  # exclude trailing paragraph mark so we don't split the paragraph
$r.Text = "Product Marketing Core $bullet Research & Analytics $bullet Communication & Technology"

# The two following paragraphs (Research & Analytics / Communication & Technology) are now
# paragraphs 7 and 8 - remove them entirely, folding their content away.
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)
$delRange = $d.Range($p7.Range.Start, $p8.Range.End)
$delRange.Delete()

# --- Part 2: Append a new TECHNICAL SKILLS section at the end of the document ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$heading = $d.Paragraphs.Item($d.Paragraphs.Count)
$heading.Range.Text = "TECHNICAL SKILLS"
$heading.Style = "Heading 2"

$heading.Range.InsertParagraphAfter()
$skills1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$skills1.Style = "Normal"
$skills1.Range.Text = "PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development; Cross-functional Team Leadership & Collaboration; Sales Enablement & Training Material Development; Data-Driven Decision Making & Analytics Interpretation"

$skills1.Range.InsertParagraphAfter()
$skills2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$skills2.Style = "Normal"
$skills2.Range.Text = "RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; Customer Journey Mapping & Behavioral Analysis; Statistical Modeling & Trend Analysis; Performance Metrics & Dashboard Development; A/B Testing & Conversion Optimization"

$skills2.Range.InsertParagraphAfter()
$skills3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$skills3.Style = "Normal"
$skills3.Range.Text = "COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Stakeholder Communication & Presentation Skills; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Marketing Technology Stack Integration; Content Strategy & Thought Leadership; Client Relationship Management & Business Development"

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
